$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The workbook is a "Estado de Cuenta" (account statement). The previous
# employee/period detail rows (16-25) are removed and replaced by a new,
# refreshed data set coming from the updated "base de datos" (per commit
# message: "Elimna EC anteriores y se agregan nuevos, se modifica base de
# datos"). Each worker now appears in two consecutive rows (periods 1803
# and 1802) grouped together, instead of being split into two blocks of
# period 1802 followed by period 1803.
# ---------------------------------------------------------------------------

# r, TipoDoc, NumDoc, Nombre, Periodo, ValorMora(F), SalarioBasico(G)
$data = @(
    @(16, "CC", "8854409",    "NILSON HERRERA PEREZ",           "1803", 36800,  957628),
    @(17, "CC", "8854409",    "NILSON HERRERA PEREZ",           "1802", 38305,  957628),
    @(18, "CC", "1032408375", "HECTOR DANIEL GARCIA ABONDANO",  "1803", 166544, 6551730),
    @(19, "CC", "1032408375", "HECTOR DANIEL GARCIA ABONDANO",  "1802", 166544, 6551730),
    @(20, "CC", "52718112",   "ANGELICA MARIA GULFO BASTIDAS",  "1803", 200000, 0),
    @(21, "CC", "52718112",   "ANGELICA MARIA GULFO BASTIDAS",  "1802", 200000, 0),
    @(22, "CC", "1047471603", "LEONARDO DANIEL ADARRAGA PINTO", "1803", 52836,  1320902),
    @(23, "CC", "1047471603", "LEONARDO DANIEL ADARRAGA PINTO", "1802", 52836,  1320902),
    @(24, "CC", "73089307",   "GERMAN ANTONIO GARZON GOMEZ",    "1803", 42401,  1060023),
    @(25, "CC", "73089307",   "GERMAN ANTONIO GARZON GOMEZ",    "1802", 42401,  1060023)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value2 = $row[1]   # B: Tipo Doc Trabajador
    $ws.Cells.Item($r, 3).Value2 = $row[2]   # C: N Doc Trabajador
    $ws.Cells.Item($r, 4).Value2 = $row[3]   # D: Nombre Trabajador
    $ws.Cells.Item($r, 5).Value2 = $row[4]   # E: Periodo Mora
    $ws.Cells.Item($r, 6).Value2 = $row[5]   # F: Valor Mora
    $ws.Cells.Item($r, 7).Value2 = $row[6]   # G: Salario Basico
}

# Refresh the column widths so the resized data (new IDs / names / amounts)
# keeps being fully visible, mirroring the bestFit columns in the sheet
# (columns B, C, E, G, H, I, J grow a bit to fit the new content).
$ws.Columns.Item(2).ColumnWidth = 17.58797
$ws.Columns.Item(3).ColumnWidth = 15.76656
$ws.Columns.Item(5).ColumnWidth = 12.58797
$ws.Columns.Item(7).ColumnWidth = 13.41828
$ws.Columns.Item(8).ColumnWidth = 18.41828
$ws.Columns.Item(9).ColumnWidth = 17.25484
$ws.Columns.Item(10).ColumnWidth = 14.085
